# repull data, push all data, mean calculation
# Updates column F ("dSF") values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    4  = -1
    10 = -5
    12 = -4
    14 = 7
    16 = -3
    17 = -3
    19 = -5
    20 = -6
    22 = -7
    24 = -5
    25 = -6
    28 = -5
    30 = -6
    31 = -9
    33 = -9
    36 = -4
    38 = -6
    39 = -2
    41 = -9
    43 = 8
    45 = -1
    51 = 0
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
